# Update the "list" sheet: asset_category_type row gains an extra enum
# value/label ('A' / "appliances").
$wb = $excel.ActiveWorkbook
$listSheet = $wb.Worksheets.Item("list")

$listSheet.Range("B1").Value = "('F', 'E', 'A')"
$listSheet.Range("C1").Value = "facility, equipment, appliances"

# Move the selection to C1, as seen in the saved file.
$listSheet.Range("C1").Select()

$wb.Save()
